$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 333-334, pushing the existing rows
# (the old 333-400 block) down to 335-402.
$ws.Rows("333:334").Insert()

# Populate the newly inserted row 333 (quality "Primera") with the
# latest weekly observation.
$ws.Range("A333").Value = 8
$ws.Range("B333").Value = "Terminal La Palmera de La Serena"
$ws.Range("C333").Value = "Coquimbo"
$ws.Range("D333").Value = 44694
$ws.Range("E333").Value = 4
$ws.Range("F333").Value = 100112017
$ws.Range("G333").Value = "Apio"
$ws.Range("H333").Value = "Americana (o)"
$ws.Range("I333").Value = "Primera"
$ws.Range("J333").Value = 2200
$ws.Range("K333").Value = 8000
$ws.Range("L333").Value = 9000
$ws.Range("M333").Value = 8500
$ws.Range("N333").Value = "$/docena de matas"
$ws.Range("O333").Value = "Provincia del Elquí"
$ws.Range("P333").Value = 1417
$ws.Range("Q333").Value = 6
$ws.Range("R333").Value = "Hortaliza"

# Populate the newly inserted row 334 (quality "Segunda") with the
# same weekly observation.
$ws.Range("A334").Value = 8
$ws.Range("B334").Value = "Terminal La Palmera de La Serena"
$ws.Range("C334").Value = "Coquimbo"
$ws.Range("D334").Value = 44694
$ws.Range("E334").Value = 4
$ws.Range("F334").Value = 100112017
$ws.Range("G334").Value = "Apio"
$ws.Range("H334").Value = "Americana (o)"
$ws.Range("I334").Value = "Segunda"
$ws.Range("J334").Value = 1480
$ws.Range("K334").Value = 6000
$ws.Range("L334").Value = 7000
$ws.Range("M334").Value = 6500
$ws.Range("N334").Value = "$/docena de matas"
$ws.Range("O334").Value = "Provincia del Elquí"
$ws.Range("P334").Value = 1083
$ws.Range("Q334").Value = 6
$ws.Range("R334").Value = "Hortaliza"
